# Weekly update: insert two new fruit/hortaliza price records for
# "Vega Monumental Concepción - Mandarina" near the top of the existing
# data block (rows 211-212), pushing the previously-existing rows
# (211-234) down by two rows (to 213-236).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 211 so everything currently at 211..234
# shifts down to 213..236, carrying its values and formatting with it.
$ws.Rows("211:212").Insert()

# --- Row 211: new Clemenuless / Primera record -------------------------
$ws.Cells.Item(211, 1).Value  = 11
$ws.Cells.Item(211, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(211, 3).Value  = "Bíobío"
$ws.Cells.Item(211, 4).Value  = (Get-Date -Year 2023 -Month 7 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(211, 5).Value  = 8
$ws.Cells.Item(211, 6).Value  = "Fruta"
$ws.Cells.Item(211, 7).Value  = 100102
$ws.Cells.Item(211, 8).Value  = "Cítricos"
$ws.Cells.Item(211, 9).Value  = 100102004
$ws.Cells.Item(211, 10).Value = "Mandarina"
$ws.Cells.Item(211, 11).Value = "Clemenuless"
$ws.Cells.Item(211, 12).Value = "Primera"
$ws.Cells.Item(211, 13).Value = 140
$ws.Cells.Item(211, 14).Value = 8000
$ws.Cells.Item(211, 15).Value = 8000
$ws.Cells.Item(211, 16).Value = 8000
$ws.Cells.Item(211, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(211, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(211, 19).Value = 800
$ws.Cells.Item(211, 20).Value = 10

# --- Row 212: new Murcott / Primera record ------------------------------
$ws.Cells.Item(212, 1).Value  = 11
$ws.Cells.Item(212, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(212, 3).Value  = "Bíobío"
$ws.Cells.Item(212, 4).Value  = (Get-Date -Year 2023 -Month 7 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(212, 5).Value  = 8
$ws.Cells.Item(212, 6).Value  = "Fruta"
$ws.Cells.Item(212, 7).Value  = 100102
$ws.Cells.Item(212, 8).Value  = "Cítricos"
$ws.Cells.Item(212, 9).Value  = 100102004
$ws.Cells.Item(212, 10).Value = "Mandarina"
$ws.Cells.Item(212, 11).Value = "Murcott"
$ws.Cells.Item(212, 12).Value = "Primera"
$ws.Cells.Item(212, 13).Value = 270
$ws.Cells.Item(212, 14).Value = 10000
$ws.Cells.Item(212, 15).Value = 11000
$ws.Cells.Item(212, 16).Value = 10556
$ws.Cells.Item(212, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(212, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(212, 19).Value = 586
$ws.Cells.Item(212, 20).Value = 18
